$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the "Lower Right Cell" references in column D (rows 5-11) from row 31 to row 33
$ws.Range("D5").Value = "A33"
$ws.Range("D6").Value = "B33"
$ws.Range("D7").Value = "C33"
$ws.Range("D8").Value = "G33"
$ws.Range("D9").Value = "H33"
$ws.Range("D10").Value = "I33"
$ws.Range("D11").Value = "J33"
